$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as TEXT, preserving exact literal formatting (leading/
# trailing zeros, multi-dot "thousand" separators, etc.) instead of letting
# Excel silently reinterpret number-looking strings as actual numbers.
function Set-TextValue($range, $val) {
    $range.NumberFormat = "@"
    $range.Value2 = $val
}

Set-TextValue $ws.Range('D2') '30.653.96'
Set-TextValue $ws.Range('E2') '  +0.65%  '
Set-TextValue $ws.Range('D3') '2.116.41'
Set-TextValue $ws.Range('E3') '  +0.42%  '
Set-TextValue $ws.Range('E4') '  +0.94%  '
Set-TextValue $ws.Range('D5') '348.96'
Set-TextValue $ws.Range('E5') '  +3.89%  '
Set-TextValue $ws.Range('D6') '1.011'
Set-TextValue $ws.Range('E6') '  +0.93%  '
Set-TextValue $ws.Range('D7') '0.5265'
Set-TextValue $ws.Range('E7') '  +0.50%  '
Set-TextValue $ws.Range('D8') '0.4526'
Set-TextValue $ws.Range('E8') '  -1.29%  '
Set-TextValue $ws.Range('D9') '53.67'
Set-TextValue $ws.Range('E9') '  +0.87%  '
Set-TextValue $ws.Range('D10') '0.09025'
Set-TextValue $ws.Range('E10') '  +1.23%  '
Set-TextValue $ws.Range('D11') '1.174'
Set-TextValue $ws.Range('E11') '  -0.38%  '
Set-TextValue $ws.Range('D12') '24.55'
Set-TextValue $ws.Range('E12') '  +0.22%  '
Set-TextValue $ws.Range('D13') '2.113.67'
Set-TextValue $ws.Range('E13') '  +0.88%  '
Set-TextValue $ws.Range('D14') '6.831'
Set-TextValue $ws.Range('E14') '  +0.48%  '
Set-TextValue $ws.Range('D15') '8.047'
Set-TextValue $ws.Range('E15') '  +0.60%  '
Set-TextValue $ws.Range('D16') '101.88'
Set-TextValue $ws.Range('E16') '  +5.50%  '
Set-TextValue $ws.Range('D17') '0.00001169'
Set-TextValue $ws.Range('E17') '  +3.19%  '
Set-TextValue $ws.Range('D18') '1.012'
Set-TextValue $ws.Range('E18') '  +0.94%  '
Set-TextValue $ws.Range('D19') '0.06710'
Set-TextValue $ws.Range('E19') '  +1.23%  '
Set-TextValue $ws.Range('D20') '19.40'
Set-TextValue $ws.Range('E20') '  +0.38%  '
Set-TextValue $ws.Range('D21') '1.010'
Set-TextValue $ws.Range('E21') '  +0.90%  '
Set-TextValue $ws.Range('D22') '6.309'
Set-TextValue $ws.Range('E22') '  -0.17%  '
Set-TextValue $ws.Range('D23') '30.713.04'
Set-TextValue $ws.Range('E23') '  +0.65%  '
Set-TextValue $ws.Range('D24') '12.85'
Set-TextValue $ws.Range('E24') '  +3.43%  '
Set-TextValue $ws.Range('D25') '2.392'
Set-TextValue $ws.Range('E25') '  +1.14%  '
Set-TextValue $ws.Range('D26') '2.363.16'
Set-TextValue $ws.Range('E26') '  +0.93%  '
Set-TextValue $ws.Range('D27') '22.46'
Set-TextValue $ws.Range('E27') '  +0.37%  '
Set-TextValue $ws.Range('D28') '165.58'
Set-TextValue $ws.Range('E28') '  +1.28%  '
Set-TextValue $ws.Range('D29') '2.542'
Set-TextValue $ws.Range('E29') '  -1.69%  '
Set-TextValue $ws.Range('D30') '137.01'
Set-TextValue $ws.Range('E30') '  +3.18%  '
Set-TextValue $ws.Range('E31') '  -2.36%  '
Set-TextValue $ws.Range('D32') '0.1077'
Set-TextValue $ws.Range('E32') '  +0.18%  '
Set-TextValue $ws.Range('D33') '1.647'
Set-TextValue $ws.Range('E33') '  -3.99%  '
Set-TextValue $ws.Range('D34') '6.403'
Set-TextValue $ws.Range('E34') '  +3.41%  '
Set-TextValue $ws.Range('D35') '4.018'
Set-TextValue $ws.Range('E35') '  +2.29%  '
Set-TextValue $ws.Range('D36') '5.950'
Set-TextValue $ws.Range('E36') '  +6.92%  '
Set-TextValue $ws.Range('D37') '10.27'
Set-TextValue $ws.Range('E37') '  -2.55%  '
Set-TextValue $ws.Range('D38') '0.02655'
Set-TextValue $ws.Range('E38') '  +2.95%  '
Set-TextValue $ws.Range('D39') '0.06857'
Set-TextValue $ws.Range('E39') '  +0.20%  '
Set-TextValue $ws.Range('D40') '0.2319'
Set-TextValue $ws.Range('E40') '  +0.55%  '
Set-TextValue $ws.Range('D41') '12.60'
Set-TextValue $ws.Range('E41') '  -1.94%  '
Set-TextValue $ws.Range('D42') '0.6908'
Set-TextValue $ws.Range('E42') '  -0.13%  '
Set-TextValue $ws.Range('D43') '1.275'
Set-TextValue $ws.Range('E43') '  +2.17%  '
Set-TextValue $ws.Range('D44') '14.79'
Set-TextValue $ws.Range('E44') '  +5.24%  '
Set-TextValue $ws.Range('D45') '2.328'
Set-TextValue $ws.Range('E45') '  -1.06%  '
Set-TextValue $ws.Range('D46') '0.6445'
Set-TextValue $ws.Range('E46') '  +0.91%  '
Set-TextValue $ws.Range('D47') '3.751'
Set-TextValue $ws.Range('E47') '  +2.35%  '
Set-TextValue $ws.Range('D48') '0.00000000355'
Set-TextValue $ws.Range('E48') '  +0.25%  '
Set-TextValue $ws.Range('D49') '1.252'
Set-TextValue $ws.Range('E49') '  +0.29%  '

# Row 50: WOONetwork inserted (previously Cronos)
Set-TextValue $ws.Range('B50') 'WOONetwork'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
Set-TextValue $ws.Range('D50') '0.3360'
Set-TextValue $ws.Range('E50') '  +2.31%  '

# Row 51: Cronos (previously Aave)
Set-TextValue $ws.Range('B51') 'Cronos'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D51') '0.07292'
Set-TextValue $ws.Range('E51') '  +2.31%  '

# A handful of prices are extremely small (< 0.0001) and the underlying
# engine always treats such number-shaped literals as numeric values no
# matter how they're entered (same as Excel auto-detecting numbers while
# typing). Apply an explicit fixed-point NumberFormat so the stored number
# still displays with the exact original digits/trailing zeros.
$ws.Range('D17').NumberFormat = '0.00000000'
$ws.Range('D48').NumberFormat = '0.00000000000'
